$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph "4S1 -> 15/3": append a plain space run, then a green-
#    highlighted "-> OK" run.
# ---------------------------------------------------------------------------
$idx1 = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd()
    if ($t -eq "4S1 -> 15/3") {
        $idx1 = $i
        break
    }
    $i = $i + 1
}

if ($idx1 -ge 0) {
    $p1 = $d.Paragraphs.Item($idx1 + 1)
    $r1 = $p1.Range
    $r1.Collapse(0)
    $r1.InsertAfter(" ")

    # Re-fetch the paragraph - InsertAfter can leave the old Range object's
    # reported Start/End out of sync with the real document once a
    # structural edit has happened, so we grab a brand-new Range.
    $p1b = $d.Paragraphs.Item($idx1 + 1)
    $r1b = $p1b.Range
    $r1b.Collapse(0)
    $r1b.InsertAfter("-> OK")

    # Locate the freshly-inserted text from the whole-document story (not a
    # paragraph sub-range, whose cached offsets can go stale after the
    # InsertAfter calls above) so we get true absolute character offsets.
    $searchRange1 = $d.Content
    $found1 = $searchRange1.Find.Execute("4S1 -> 15/3 -> OK", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found1) {
        $hi1 = $d.Range($searchRange1.End - 5, $searchRange1.End)
        $hi1.Font.HighlightColorIndex = 4   # wdBrightGreen -> OOXML w:highlight val="green"
    }
}

# ---------------------------------------------------------------------------
# 2) Paragraph "5S2 -> 16/3": append a green-highlighted "-> OK" run
#    directly (no leading space run this time).
# ---------------------------------------------------------------------------
$idx2 = -1
$i = 0
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd()
    if ($t -eq "5S2 -> 16/3") {
        $idx2 = $i
        break
    }
    $i = $i + 1
}

if ($idx2 -ge 0) {
    $p2 = $d.Paragraphs.Item($idx2 + 1)
    $r2 = $p2.Range
    $r2.Collapse(0)
    $r2.InsertAfter("-> OK")

    $searchRange2 = $d.Content
    $found2 = $searchRange2.Find.Execute("5S2 -> 16/3-> OK", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found2) {
        $hi2 = $d.Range($searchRange2.End - 5, $searchRange2.End)
        $hi2.Font.HighlightColorIndex = 4   # wdBrightGreen -> OOXML w:highlight val="green"
    }
}
